# Translations.xlsx — add translations for previously-untranslated strings.
#
# 1) A new row is inserted at row 85 for "AddMovingPlatformWaypoint"
#    (shifting the old rows 85-131 down to 86-132).
# 2) Thirteen new rows are appended after the old last row (now 132),
#    i.e. rows 133-145, for the moving-platform / door related strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert the new "AddMovingPlatformWaypoint" row at position 85 ---
$ws.Rows("85:85").Insert()

$ws.Range("A85").Value = "AddMovingPlatformWaypoint"
$ws.Range("B85").Value = "+ Add Platform Waypoint"
$ws.Range("C85").Value = "+ Añadir Waypoint de Plataforma"

# --- 2) Append the new rows 133-145 ---
$ws.Range("A133").Value = "InvisibleEdges"
$ws.Range("B133").Value = "Invisible Edges"
$ws.Range("C133").Value = "Border Invisibles"

$ws.Range("A134").Value = "StartMovingAtStart"
$ws.Range("B134").Value = "Start Moving At Start"
$ws.Range("C134").Value = "Mover Al Inicio"

$ws.Range("A135").Value = "MovingSpeed"
$ws.Range("B135").Value = "Moving Speed"
$ws.Range("C135").Value = "Velocidad de Movimiento"

$ws.Range("A136").Value = "StartDelay"
$ws.Range("B136").Value = "Start Delay"
$ws.Range("C136").Value = "Espera al Inicio"

$ws.Range("A137").Value = "MovementMode"
$ws.Range("B137").Value = "Movement Mode"
$ws.Range("C137").Value = "Modo de Movimiento"

$ws.Range("A138").Value = "None_Mayus"
$ws.Range("B138").Value = "NONE"
$ws.Range("C138").Value = "NINGUNO"

$ws.Range("A139").Value = "TravelBack_Mayus"
$ws.Range("B139").Value = "TRAVEL BACK"
$ws.Range("C139").Value = "REGRESARSE"

$ws.Range("A140").Value = "Loop_Mayus"
$ws.Range("B140").Value = "LOOP"
$ws.Range("C140").Value = "BUCLE"

$ws.Range("A141").Value = "IsAutomatic"
$ws.Range("B141").Value = "Is Automatic?"
$ws.Range("C141").Value = "¿Es Automática?"

$ws.Range("A142").Value = "CLOSED"
$ws.Range("B142").Value = "CLOSED"
$ws.Range("C142").Value = "CERRADA"

$ws.Range("A143").Value = "OPEN"
$ws.Range("B143").Value = "OPEN"
$ws.Range("C143").Value = "ABIERTA"

$ws.Range("A144").Value = "LOCKED"
$ws.Range("B144").Value = "LOCKED"
$ws.Range("C144").Value = "BLOQUEADA"

$ws.Range("A145").Value = "UNLOCKED"
$ws.Range("B145").Value = "UNLOCKED"
$ws.Range("C145").Value = "DESBLOQUEADA"

# --- 3) Update the view state to match what Excel would save after this
#        editing session: selection resting on the next empty row, and the
#        window scrolled down so row 126 is the first visible row. ---
$ws.Activate()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 126 } catch {}
try { $win.ScrollColumn = 1 } catch {}
[void]$ws.Range("A146").Select()
